$wb = $excel.ActiveWorkbook

# --- Sheets ---
$metadata = $wb.Worksheets.Item("Metadata")
$core = $wb.Worksheets.Item("Core")
$chartOne = $wb.Worksheets.Item("Test Chart One")
$chartTwo = $wb.Worksheets.Item("Test Chart Two")

$xlPasteFormats = -4122

# --- Style changes (direct cell formats) ---
# This workbook only has two distinct direct cell formats in use ("plain" and
# an "alternate" one with an explicit alignment). We reuse the existing
# format definitions by copy/paste-special instead of re-deriving them, and
# we must grab the "alternate" style from still-untouched source cells
# *before* those very cells get reformatted below.

# 1) Capture the alternate style from Metadata!A8 (still alternate at this
#    point) and apply it to the two cells that need to gain that style.
$metadata.Range("A8").Copy() | Out-Null
$chartOne.Range("A2").PasteSpecial($xlPasteFormats) | Out-Null
$chartTwo.Range("A2").PasteSpecial($xlPasteFormats) | Out-Null

# 2) Now move every cell that should lose the alternate style back to the
#    plain style, copied from a cell that keeps the plain style throughout.
$metadata.Range("A7").Copy() | Out-Null
$metadata.Range("A8:B9").PasteSpecial($xlPasteFormats) | Out-Null
$core.Range("R5").PasteSpecial($xlPasteFormats) | Out-Null
$chartTwo.Range("C3:D4").PasteSpecial($xlPasteFormats) | Out-Null

$excel.CutCopyMode = 0

# --- Value changes ---
# Core sheet: column A should hold the actual question code (matching column
# B) instead of the old hard-coded test id.
$core.Cells.Item(2, 1).Value = "ComplexChartInstanceName"
$core.Cells.Item(3, 1).Value = "ComplexChartDate"
$core.Cells.Item(4, 1).Value = "ComplexChartType"
$core.Cells.Item(5, 1).Value = "ComplexChartSubtype"

# Charting date question now uses a constant id across both chart sheets.
$chartOne.Cells.Item(2, 1).Value = "PatientChartingDate"
$chartTwo.Cells.Item(2, 1).Value = "PatientChartingDate"
